$d = $word.ActiveDocument

# 1) "Краткое описание" -> "КРАТКОЕ ОПИСАНИЕ"
$d.Content.Find.Execute("Краткое описание", $true, $false, $false, $false, $false,
                         $true, 1, $false, "КРАТКОЕ ОПИСАНИЕ", 2) | Out-Null

# 2) "Основной поток событий" -> "ОСНОВНОЙ ПОТОК СОБЫТИЙ"
$d.Content.Find.Execute("Основной поток событий", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ОСНОВНОЙ ПОТОК СОБЫТИЙ", 2) | Out-Null

# 3) "Альтернативные потоки" (spread across 4 runs) -> single run "АЛЬТЕРНАТИВНЫЕ ПОТОКИ"
$d.Content.Find.Execute("Альтернативные потоки", $true, $false, $false, $false, $false,
                         $true, 1, $false, "АЛЬТЕРНАТИВНЫЕ ПОТОКИ", 2) | Out-Null

# 4) "Предусловия" -> "ПРЕДУСЛОВИЯ"
$d.Content.Find.Execute("Предусловия", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ПРЕДУСЛОВИЯ", 2) | Out-Null

# 5) "Постусловия" -> "ПОСТУСЛОВИЯ"
$d.Content.Find.Execute("Постусловия", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ПОСТУСЛОВИЯ", 2) | Out-Null

# 6) Move the "_GoBack" bookmark from the end of the final paragraph to right
#    after the "ПОСТУСЛОВИЯ" run (end of that paragraph's text).
#    Bookmarks.Add requires a non-collapsed Range, so we temporarily insert a
#    marker character at the end of the "ПОСТУСЛОВИЯ" paragraph, bookmark that
#    single character (keeping the "ПОСТУСЛОВИЯ" run untouched), then delete
#    the marker again. Re-adding "_GoBack" relocates the existing bookmark
#    (Word only ever keeps a single "_GoBack") instead of creating a duplicate.
$post = $d.Paragraphs(18).Range
$insPt = $d.Range($post.End - 1, $post.End - 1)
$insPt.InsertAfter("X")

$post2 = $d.Paragraphs(18).Range
$markerStart = $post2.End - 2
$markerEnd = $post2.End - 1
$marker = $d.Range($markerStart, $markerEnd)
$d.Bookmarks.Add("_GoBack", $marker) | Out-Null

$d.Range($markerStart, $markerEnd).Delete() | Out-Null

Write-Output "done"
